$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I15").Value = 130.38
$ws.Range("H15").Value = 130.38
$ws.Range("M15").Value = -222.14
$ws.Range("K15").Value = 391.14
$ws.Range("N17").Value = -5962.9353
$ws.Range("L17").Value = 5626.9353
$ws.Range("H17").Value = 1875.6451
$ws.Range("J17").Value = 1875.6451
$ws.Range("K40").Value = 1592.8572
$ws.Range("N40").Value = -2850
$ws.Range("L40").Value = 2500
$ws.Range("H40").Value = 1922.7273
$ws.Range("J40").Value = 2500
$ws.Range("M40").Value = -1417.8572
$ws.Range("I40").Value = 1592.8572
$ws.Range("K107").Value = 111653.89
$ws.Range("H107").Value = 91401
$ws.Range("M107").Value = -109733.89
$ws.Range("I107").Value = 111653.89
$ws.Range("K113").Value = 1458
$ws.Range("M113").Value = 1796
$ws.Range("I113").Value = 1458
$ws.Range("H113").Value = 2571.818
$ws.Range("N127").Value = -14187.4517
$ws.Range("L127").Value = 4267.4517
$ws.Range("J127").Value = 1422.4839
$ws.Range("H127").Value = 1354.4546
$ws.Range("N129").Value = -13391.1334
$ws.Range("L129").Value = 3391.1334
$ws.Range("J129").Value = 1130.3778
$ws.Range("H129").Value = 1120.5745
$ws.Range("I132").Value = 1700.5385
$ws.Range("H132").Value = 1837.138
$ws.Range("K132").Value = 5101.6155
$ws.Range("M132").Value = -2571.6155
$ws.Range("N134").Value = -129691.97
$ws.Range("L134").Value = 119551.97
$ws.Range("J134").Value = 119551.97
$ws.Range("H134").Value = 119551.97
$ws.Range("K137").Value = 4956.3915
$ws.Range("N137").Value = -14398.8
$ws.Range("H137").Value = 1910.6072
$ws.Range("L137").Value = 9298.799999999999
$ws.Range("J137").Value = 3099.6
$ws.Range("M137").Value = -2406.3915
$ws.Range("I137").Value = 1652.1305
$ws.Range("H138").Value = 11120359
$ws.Range("J138").Value = 10746.786
$ws.Range("M138").Value = -150006860
$ws.Range("I138").Value = 50004000
$ws.Range("K138").Value = 150012000
$ws.Range("N138").Value = -42520.358
$ws.Range("L138").Value = 32240.358
$ws.Range("M141").Value = -10853.3329
$ws.Range("I141").Value = 5344.4443
$ws.Range("H141").Value = 10310
$ws.Range("K141").Value = 16033.3329

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J74").Value = 1137.3334
$ws.Range("M74").Value = -213.075
$ws.Range("I74").Value = 1087.075
$ws.Range("K74").Value = 1087.075
$ws.Range("N74").Value = -2885.3334
$ws.Range("L74").Value = 1137.3334
$ws.Range("H74").Value = 1093.6305
$ws.Range("M77").Value = -1067.375
$ws.Range("I77").Value = 1087.075
$ws.Range("H77").Value = 1093.6305
$ws.Range("K77").Value = 5435.375
$ws.Range("N77").Value = -14422.667
$ws.Range("L77").Value = 5686.666999999999
$ws.Range("J77").Value = 1137.3334
$ws.Range("N123").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("I132").Value = 1288.9166
$ws.Range("H132").Value = 1442.1642
$ws.Range("K132").Value = 3866.7498
$ws.Range("M132").Value = -1336.7498
$ws.Range("N134").ClearContents()
$ws.Range("L134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("H134").Value = 0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N53").ClearContents()
$ws.Range("L53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("H64").Value = 238
$ws.Range("M64").Value = -61.58334000000002
$ws.Range("I64").Value = 286.58334
$ws.Range("K64").Value = 286.58334
$ws.Range("H67").Value = 238
$ws.Range("M67").Value = 493.41666
$ws.Range("I67").Value = 286.58334
$ws.Range("K67").Value = 286.58334

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N31").Value = -6630.25
$ws.Range("L31").Value = 6040.25
$ws.Range("H31").Value = 4398.524
$ws.Range("J31").Value = 6040.25
$ws.Range("M31").Value = -3093.2307
$ws.Range("I31").Value = 3388.2307
$ws.Range("K31").Value = 3388.2307
$ws.Range("K34").Value = 3388.2307
$ws.Range("N34").Value = -6444.25
$ws.Range("L34").Value = 6040.25
$ws.Range("I34").Value = 3388.2307
$ws.Range("H34").Value = 4398.524
$ws.Range("J34").Value = 6040.25
$ws.Range("M34").Value = -3186.2307
$ws.Range("L70").Value = 12000
$ws.Range("J70").Value = 12000
$ws.Range("H70").Value = 12000
$ws.Range("N70").Value = -12630
$ws.Range("N73").Value = -14184
$ws.Range("H73").Value = 12000
$ws.Range("L73").Value = 12000
$ws.Range("J73").Value = 12000
$ws.Range("I132").Value = 1872.44
$ws.Range("H132").Value = 2079.4285
$ws.Range("K132").Value = 5617.32
$ws.Range("N132").Value = -16472.9999
$ws.Range("L132").Value = 11412.9999
$ws.Range("J132").Value = 3804.3333
$ws.Range("M132").Value = -3087.32
$ws.Range("N134").Value = -10855.5
$ws.Range("L134").Value = 5785.5
$ws.Range("J134").Value = 1928.5
$ws.Range("M134").Value = -2510.3748
$ws.Range("I134").Value = 1681.7916
$ws.Range("H134").Value = 1717.0358
$ws.Range("K134").Value = 5045.3748

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 2867.6667
$ws.Range("J48").Value = 3500
$ws.Range("M48").Value = -4559
$ws.Range("I48").Value = 1603
$ws.Range("K48").Value = 4809
$ws.Range("N48").Value = -11000
$ws.Range("L48").Value = 10500
$ws.Range("N75").Value = -36996.001
$ws.Range("L75").Value = 35000.001
$ws.Range("J75").Value = 11666.667
$ws.Range("H75").Value = 5487
$ws.Range("N78").Value = -114984.003
$ws.Range("H78").Value = 5487
$ws.Range("L78").Value = 105000.003
$ws.Range("J78").Value = 11666.667
$ws.Range("K113").Value = 2600.0001
$ws.Range("N113").Value = -60739.50199999999
$ws.Range("L113").Value = 56399.50199999999
$ws.Range("J113").Value = 18799.834
$ws.Range("M113").Value = -430.0001000000002
$ws.Range("I113").Value = 866.6667
$ws.Range("H113").Value = 8039.933
$ws.Range("L131").Value = 37503795
$ws.Range("H131").Value = 12348163
$ws.Range("J131").Value = 12501265
$ws.Range("N131").Value = -37513875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N15").Value = -264076
$ws.Range("L15").Value = 263500
$ws.Range("H15").Value = 263500
$ws.Range("J15").Value = 263500
$ws.Range("N81").Value = -265496
$ws.Range("L81").Value = 263500
$ws.Range("H81").Value = 263500
$ws.Range("J81").Value = 263500
$ws.Range("N84").Value = -800484
$ws.Range("L84").Value = 790500
$ws.Range("J84").Value = 263500
$ws.Range("H84").Value = 263500
$ws.Range("M122").Value = -550
$ws.Range("I122").Value = 1000
$ws.Range("H122").Value = 4801.6
$ws.Range("K122").Value = 3000
$ws.Range("N122").Value = -22156
$ws.Range("L122").Value = 17256
$ws.Range("J122").Value = 5752
$ws.Range("N123").Value = -22230.646
$ws.Range("H123").Value = 17330.646
$ws.Range("L123").Value = 17330.646
$ws.Range("J123").Value = 17330.646

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J16").Value = 444
$ws.Range("M16").Value = -1442.3334
$ws.Range("I16").Value = 1612.3334
$ws.Range("K16").Value = 1612.3334
$ws.Range("H16").Value = 1195.0714
$ws.Range("N16").Value = -784
$ws.Range("L16").Value = 444
$ws.Range("K40").Value = 10829.583
$ws.Range("N40").Value = -10247.556
$ws.Range("L40").Value = 9975.556
$ws.Range("H40").Value = 10463.571
$ws.Range("J40").Value = 9975.556
$ws.Range("M40").Value = -10693.583
$ws.Range("I40").Value = 10829.583
$ws.Range("K46").Value = 1098.6
$ws.Range("N46").Value = -1990.5714
$ws.Range("L46").Value = 1614.5714
$ws.Range("J46").Value = 1614.5714
$ws.Range("I46").Value = 1098.6
$ws.Range("H46").Value = 1399.5834
$ws.Range("M46").Value = -910.5999999999999
$ws.Range("H80").Value = 34000
$ws.Range("L80").Value = 34000
$ws.Range("J80").Value = 34000
$ws.Range("N80").Value = -36246
$ws.Range("K82").Value = 1387.7778
$ws.Range("N82").Value = -2861.9
$ws.Range("L82").Value = 2139.9
$ws.Range("J82").Value = 2139.9
$ws.Range("M82").Value = -1026.7778
$ws.Range("I82").Value = 1387.7778
$ws.Range("H82").Value = 1783.6316
$ws.Range("N83").Value = -113232
$ws.Range("H83").Value = 34000
$ws.Range("L83").Value = 102000
$ws.Range("J83").Value = 34000
$ws.Range("K85").Value = 1387.7778
$ws.Range("N85").Value = -4635.9
$ws.Range("L85").Value = 2139.9
$ws.Range("H85").Value = 1783.6316
$ws.Range("J85").Value = 2139.9
$ws.Range("M85").Value = -139.7778000000001
$ws.Range("I85").Value = 1387.7778
$ws.Range("I132").Value = 4468.976
$ws.Range("H132").Value = 4482.8936
$ws.Range("K132").Value = 13406.928
$ws.Range("N132").Value = -18859.4
$ws.Range("L132").Value = 13799.4
$ws.Range("J132").Value = 4599.8
$ws.Range("M132").Value = -10876.928
$ws.Range("K136").Value = 18922.1061
$ws.Range("I136").Value = 6307.3687
$ws.Range("H136").Value = 6278.095
$ws.Range("M136").Value = -16372.1061
$ws.Range("N140").Value = -79943.336
$ws.Range("L140").Value = 69583.336
$ws.Range("H140").Value = 69583.336
$ws.Range("J140").Value = 69583.336

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K62").Value = 4000
$ws.Range("N62").Value = -4447.8
$ws.Range("L62").Value = 3199.8
$ws.Range("H62").Value = 3428.4285
$ws.Range("J62").Value = 3199.8
$ws.Range("M62").Value = -3376
$ws.Range("I62").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("N65").Value = -22239
$ws.Range("L65").Value = 15999
$ws.Range("J65").Value = 3199.8
$ws.Range("I65").Value = 4000
$ws.Range("H65").Value = 3428.4285
$ws.Range("M65").Value = -16880
$ws.Range("K123").Value = 20390
$ws.Range("H123").Value = 40259.082
$ws.Range("M123").Value = -15490
$ws.Range("I123").Value = 20390
$ws.Range("I132").Value = 1395.1875
$ws.Range("H132").Value = 1567.1428
$ws.Range("K132").Value = 4555.0344
$ws.Range("N132").Value = -11412.2
$ws.Range("L132").Value = 6352.200000000001
$ws.Range("J132").Value = 2117.4
$ws.Range("M132").Value = -1655.5625
